$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the existing
# "Late" / "heading" / "Outstanding" columns one position to the right.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 9.83

# The "Repayment schedule" tab becomes the active tab (was "Transactions"),
# with a new selected cell.
$ws.Activate()
$ws.Range("L14").Select() | Out-Null
